# Actualización automática 2025-08-04 17:26:10
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Column widths (cols D, E, F). The ColumnWidth COM setter round-trips to the
# saved XML `width` attribute with a constant +5/6 offset, so subtract it here
# to land on the exact target widths (11, 22, 18) in the saved file.
$ws.Columns.Item(4).ColumnWidth = 11 - (5/6)
$ws.Columns.Item(5).ColumnWidth = 22 - (5/6)
$ws.Columns.Item(6).ColumnWidth = 18 - (5/6)

# Row 2: 240X120 PORCELANATO - presupuesto (C) unchanged, venta (D) -> 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 344.284604629486
$ws.Range("F2").Value = 0

# Row 3: 240X80 PORCELANATO - presupuesto (C) changes
$ws.Range("C3").Value = 7120.1145
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 7120.1145
$ws.Range("F3").Value = 0

# Row 12: PANELES DECORATIVOS
$ws.Range("C12").Value = 100
$ws.Range("E12").Value = 100

# Row 13: PANELES PU
$ws.Range("C13").Value = 20
$ws.Range("E13").Value = 20

# Row 14: PANELES PVC
$ws.Range("C14").Value = 100
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 100
$ws.Range("F14").Value = 0

# Row 15: PIEDRA SINTERIZADA
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 1638
$ws.Range("F15").Value = 0

# Row 16: PORCELANATO
$ws.Range("C16").Value = 23904.58
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 23904.58
$ws.Range("F16").Value = 0

# Row 17: PUERTAS DE SEGURIDAD
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 342
$ws.Range("F17").Value = 0

# Row 19: TOTAL
$ws.Range("C19").Value = 37500.00093005039
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 37500.00093005039
$ws.Range("F19").Value = 0
